$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3464964993005633
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 157.8057217802531
$ws.Range("E2").Value = 198602002.3250627
$ws.Range("G2").Value = 198602162.1305174
